$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (B1:E1): new column headers, unbolded style copied from
# the row-label style used in column A (fontId 2 / style index 3) ---
$ws.Range("B1").Value = "Revenue (£)"
$ws.Range("C1").Value = "Cost of Sales (CoS) (£)"
$ws.Range("D1").Value = "Gross Profit (£)"
$ws.Range("E1").Value = "Gross Margin (%)"
$ws.Range("A2").Copy()
$ws.Range("B1:E1").PasteSpecial(-4122)  # xlPasteFormats

# --- Row labels (A2:A13): months instead of metric names ---
$ws.Range("A2").Value = "Jan"
$ws.Range("A3").Value = "Feb"
$ws.Range("A4").Value = "Mar"
$ws.Range("A5").Value = "Apr"
$ws.Range("A6").Value = "May"
$ws.Range("A7").Value = "Jun"
$ws.Range("A8").Value = "Jul"
$ws.Range("A9").Value = "Aug"
$ws.Range("A10").Value = "Sep"
$ws.Range("A11").Value = "Oct"
$ws.Range("A12").Value = "Nov"
$ws.Range("A13").Value = "Dec"

# --- Drop the now-unused extra metric rows ---
$ws.Rows("14:19").Delete()

# --- Column widths for the newly widened header columns (best-fit to content) ---
$ws.Columns("B").ColumnWidth = 9.3
$ws.Columns("C").ColumnWidth = 16.6
$ws.Columns("D").ColumnWidth = 11.5

# --- Selection moves to F27 ---
$ws.Range("F27").Select() | Out-Null
